$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38 gets new data values (was all zeros for B:N), matching merged data update 2020_04_06
$ws.Range("B38").Value = 43927
$ws.Range("B2").Copy()
$ws.Range("B38").PasteSpecial(-4122)

$ws.Range("C38").Value = 1345048
$ws.Range("D38").Value = 276515
$ws.Range("E38").Value = 74565
$ws.Range("F38").Value = 366614
$ws.Range("G38").Value = 19581
$ws.Range("H38").Value = 10783
$ws.Range("I38").Value = 131815
$ws.Range("K38").Value = 4698
$ws.Range("L38").Value = 4875
$ws.Range("M38").Value = 1489
$ws.Range("N38").Value = 187
